# Auto-generated edit script: updates currentAveragePrice-derived
# market/profit columns (H-N) for the rows touched by the scheduled
# market-data refresh across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 59268780
$ws.Cells.Item(86, 9).Value = 64010000
$ws.Cells.Item(86, 11).Value = 64010000
$ws.Cells.Item(86, 13).Value = -64008877
$ws.Cells.Item(89, 8).Value = 59268780
$ws.Cells.Item(89, 9).Value = 64010000
$ws.Cells.Item(89, 11).Value = 320050000
$ws.Cells.Item(89, 13).Value = -320044384
$ws.Cells.Item(98, 8).Value = 729.5789
$ws.Cells.Item(98, 9).Value = 624.2
$ws.Cells.Item(98, 10).Value = 1124.75
$ws.Cells.Item(98, 11).Value = 624.2
$ws.Cells.Item(98, 12).Value = 1124.75
$ws.Cells.Item(98, 13).Value = 873.8
$ws.Cells.Item(98, 14).Value = -4120.75
$ws.Cells.Item(122, 8).Value = 729.5789
$ws.Cells.Item(122, 9).Value = 624.2
$ws.Cells.Item(122, 10).Value = 1124.75
$ws.Cells.Item(122, 11).Value = 1872.6
$ws.Cells.Item(122, 12).Value = 3374.25
$ws.Cells.Item(122, 13).Value = 577.3999999999999
$ws.Cells.Item(122, 14).Value = -8274.25
$ws.Cells.Item(138, 8).Value = 2858.3828
$ws.Cells.Item(138, 9).Value = 1520.0294
$ws.Cells.Item(138, 11).Value = 4560.0882
$ws.Cells.Item(138, 13).Value = 579.9117999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(39, 8).Value = 13005.333
$ws.Cells.Item(39, 9).Value = 13005.333
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 13005.333
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).Value = ""
$ws.Cells.Item(39, 14).Value = -12485.333
$ws.Cells.Item(97, 8).Value = 1193.3334
$ws.Cells.Item(97, 9).Value = 1095.2
$ws.Cells.Item(97, 10).Value = 1473.7142
$ws.Cells.Item(97, 11).Value = 1095.2
$ws.Cells.Item(97, 12).Value = 1473.7142
$ws.Cells.Item(97, 13).Value = -599.2
$ws.Cells.Item(97, 14).Value = -2465.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 880.55554
$ws.Cells.Item(94, 9).Value = 846.2381
$ws.Cells.Item(94, 10).Value = 928.6
$ws.Cells.Item(94, 11).Value = 846.2381
$ws.Cells.Item(94, 12).Value = 928.6
$ws.Cells.Item(94, 13).Value = -395.2381
$ws.Cells.Item(94, 14).Value = -1830.6
$ws.Cells.Item(99, 8).Value = 1618.4615
$ws.Cells.Item(99, 9).Value = 1387.6666
$ws.Cells.Item(99, 10).Value = 2137.75
$ws.Cells.Item(99, 11).Value = 1387.6666
$ws.Cells.Item(99, 12).Value = 2137.75
$ws.Cells.Item(99, 13).Value = 110.3334
$ws.Cells.Item(99, 14).Value = -5133.75
$ws.Cells.Item(103, 8).Value = 50153.75
$ws.Cells.Item(103, 10).Value = 50153.75
$ws.Cells.Item(103, 12).Value = 50153.75
$ws.Cells.Item(103, 14).Value = -52497.75
$ws.Cells.Item(134, 8).Value = 1963.1311
$ws.Cells.Item(134, 9).Value = 1495.5135
$ws.Cells.Item(134, 11).Value = 4486.5405
$ws.Cells.Item(134, 13).Value = -1951.5405

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6511.204
$ws.Cells.Item(31, 9).Value = 1215.7391
$ws.Cells.Item(31, 10).Value = 11195.654
$ws.Cells.Item(31, 11).Value = 1215.7391
$ws.Cells.Item(31, 12).Value = 11195.654
$ws.Cells.Item(31, 13).Value = -920.7391
$ws.Cells.Item(31, 14).Value = -11785.654
$ws.Cells.Item(34, 8).Value = 6511.204
$ws.Cells.Item(34, 9).Value = 1215.7391
$ws.Cells.Item(34, 10).Value = 11195.654
$ws.Cells.Item(34, 11).Value = 1215.7391
$ws.Cells.Item(34, 12).Value = 11195.654
$ws.Cells.Item(34, 13).Value = -1013.7391
$ws.Cells.Item(34, 14).Value = -11599.654
$ws.Cells.Item(124, 8).Value = 29576.5
$ws.Cells.Item(124, 10).Value = 29576.5
$ws.Cells.Item(124, 12).Value = 29576.5
$ws.Cells.Item(124, 14).Value = -34486.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 567
$ws.Cells.Item(13, 9).Value = 453.16666
$ws.Cells.Item(13, 10).Value = 680.8333
$ws.Cells.Item(13, 11).Value = 1359.49998
$ws.Cells.Item(13, 12).Value = 2042.4999
$ws.Cells.Item(13, 13).Value = -1191.49998
$ws.Cells.Item(13, 14).Value = -2378.4999
$ws.Cells.Item(82, 8).Value = 3000
$ws.Cells.Item(82, 10).Value = 5000
$ws.Cells.Item(82, 12).Value = 15000
$ws.Cells.Item(82, 14).Value = -15812
$ws.Cells.Item(85, 8).Value = 3000
$ws.Cells.Item(85, 10).Value = 5000
$ws.Cells.Item(85, 12).Value = 15000
$ws.Cells.Item(85, 14).Value = -17808
$ws.Cells.Item(117, 8).Value = 345.8
$ws.Cells.Item(117, 9).Value = 243
$ws.Cells.Item(117, 10).Value = 500
$ws.Cells.Item(117, 11).Value = 729
$ws.Cells.Item(117, 12).Value = 1500
$ws.Cells.Item(117, 13).Value = 2713
$ws.Cells.Item(117, 14).Value = -8384
$ws.Cells.Item(121, 8).Value = 1250.6296
$ws.Cells.Item(121, 9).Value = 479.85715
$ws.Cells.Item(121, 10).Value = 1520.4
$ws.Cells.Item(121, 11).Value = 1439.57145
$ws.Cells.Item(121, 12).Value = 4561.200000000001
$ws.Cells.Item(121, 13).Value = -129.5714499999999
$ws.Cells.Item(121, 14).Value = -7181.200000000001
$ws.Cells.Item(125, 8).Value = 2674.1428
$ws.Cells.Item(125, 9).Value = 1000
$ws.Cells.Item(125, 10).Value = 2953.1667
$ws.Cells.Item(125, 11).Value = 3000
$ws.Cells.Item(125, 12).Value = 8859.500100000001
$ws.Cells.Item(125, 13).Value = 1920
$ws.Cells.Item(125, 14).Value = -18699.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 420.36365
$ws.Cells.Item(2, 9).Value = 56.142857
$ws.Cells.Item(2, 10).Value = 1057.75
$ws.Cells.Item(2, 11).Value = 56.142857
$ws.Cells.Item(2, 12).Value = 1057.75
$ws.Cells.Item(2, 13).Value = 56.857143
$ws.Cells.Item(2, 14).Value = -1283.75
$ws.Cells.Item(80, 8).Value = 31864406
$ws.Cells.Item(80, 9).Value = 63626440
$ws.Cells.Item(80, 10).Value = 102372
$ws.Cells.Item(80, 11).Value = 63626440
$ws.Cells.Item(80, 12).Value = 102372
$ws.Cells.Item(80, 13).Value = -63625442
$ws.Cells.Item(80, 14).Value = -104368
$ws.Cells.Item(83, 8).Value = 31864406
$ws.Cells.Item(83, 9).Value = 63626440
$ws.Cells.Item(83, 10).Value = 102372
$ws.Cells.Item(83, 11).Value = 318132200
$ws.Cells.Item(83, 12).Value = 511860
$ws.Cells.Item(83, 13).Value = -318127208
$ws.Cells.Item(83, 14).Value = -521844
$ws.Cells.Item(105, 8).Value = 27500
$ws.Cells.Item(105, 9).Value = 10000
$ws.Cells.Item(105, 11).Value = 10000
$ws.Cells.Item(105, 13).Value = -6506
$ws.Cells.Item(111, 8).Value = 25000
$ws.Cells.Item(111, 10).Value = 25000
$ws.Cells.Item(111, 12).Value = 25000
$ws.Cells.Item(111, 14).Value = -31134
$ws.Cells.Item(119, 8).Value = 59800
$ws.Cells.Item(119, 10).Value = 59800
$ws.Cells.Item(119, 12).Value = 59800
$ws.Cells.Item(119, 14).Value = -69476
$ws.Cells.Item(132, 8).Value = 3392.9524
$ws.Cells.Item(132, 9).Value = 3308.4285
$ws.Cells.Item(132, 11).Value = 9925.2855
$ws.Cells.Item(132, 13).Value = -7395.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(4, 8).Value = 51755
$ws.Cells.Item(4, 9).Value = 23500
$ws.Cells.Item(4, 11).Value = 23500
$ws.Cells.Item(4, 13).Value = -23387
$ws.Cells.Item(28, 8).Value = 51755
$ws.Cells.Item(28, 9).Value = 23500
$ws.Cells.Item(28, 11).Value = 23500
$ws.Cells.Item(28, 13).Value = -23268
$ws.Cells.Item(37, 8).Value = 51755
$ws.Cells.Item(37, 9).Value = 23500
$ws.Cells.Item(37, 11).Value = 23500
$ws.Cells.Item(37, 13).Value = -23393
$ws.Cells.Item(82, 8).Value = 33334972
$ws.Cells.Item(82, 9).Value = 50001660
$ws.Cells.Item(82, 10).Value = 1600
$ws.Cells.Item(82, 11).Value = 50001660
$ws.Cells.Item(82, 12).Value = 1600
$ws.Cells.Item(82, 13).Value = -50001299
$ws.Cells.Item(82, 14).Value = -2322
$ws.Cells.Item(85, 8).Value = 33334972
$ws.Cells.Item(85, 9).Value = 50001660
$ws.Cells.Item(85, 10).Value = 1600
$ws.Cells.Item(85, 11).Value = 50001660
$ws.Cells.Item(85, 12).Value = 1600
$ws.Cells.Item(85, 13).Value = -50000412
$ws.Cells.Item(85, 14).Value = -4096
$ws.Cells.Item(93, 8).Value = 52000
$ws.Cells.Item(93, 9).Value = 100000
$ws.Cells.Item(93, 10).Value = 4000
$ws.Cells.Item(93, 11).Value = 100000
$ws.Cells.Item(93, 12).Value = 4000
$ws.Cells.Item(93, 13).Value = -98752
$ws.Cells.Item(93, 14).Value = -6496

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3647.35
$ws.Cells.Item(81, 9).Value = 3653.125
$ws.Cells.Item(81, 10).Value = 3624.25
$ws.Cells.Item(81, 11).Value = 7306.25
$ws.Cells.Item(81, 12).Value = 7248.5
$ws.Cells.Item(81, 13).Value = -6245.25
$ws.Cells.Item(81, 14).Value = -9370.5
$ws.Cells.Item(84, 8).Value = 3647.35
$ws.Cells.Item(84, 9).Value = 3653.125
$ws.Cells.Item(84, 10).Value = 3624.25
$ws.Cells.Item(84, 11).Value = 36531.25
$ws.Cells.Item(84, 12).Value = 36242.5
$ws.Cells.Item(84, 13).Value = -31227.25
$ws.Cells.Item(84, 14).Value = -46850.5
$ws.Cells.Item(116, 8).Value = 80000
$ws.Cells.Item(116, 10).Value = 80000
$ws.Cells.Item(116, 12).Value = 80000
$ws.Cells.Item(116, 14).Value = -89178
$ws.Cells.Item(132, 8).Value = 2733989.2
$ws.Cells.Item(132, 9).Value = 1767.2894
$ws.Cells.Item(132, 11).Value = 5301.8682
$ws.Cells.Item(132, 13).Value = -2771.8682
$ws.Cells.Item(136, 8).Value = 2110.614
$ws.Cells.Item(136, 9).Value = 1652.6171
$ws.Cells.Item(136, 10).Value = 4263.2
$ws.Cells.Item(136, 11).Value = 4957.8513
$ws.Cells.Item(136, 12).Value = 12789.6
$ws.Cells.Item(136, 13).Value = -2407.8513
$ws.Cells.Item(136, 14).Value = -17889.6

